$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APRIL")

# --- Fill in the previously-missing "Current Environment" / "Who Detected"
#     columns for the existing defect row (row 6) ---
$ws.Cells.Item(6,8).Value = "Android ICS"
$ws.Cells.Item(6,9).Value = "Small"

# --- Add a new defect entry as row 7 ---
# Copy the date cell's formatting (style) from the row above so the new
# date cell renders the same way (m/d/yyyy) as the rest of the sheet.
$ws.Cells.Item(6,3).Copy()
$ws.Cells.Item(7,3).PasteSpecial(-4122)

$ws.Cells.Item(7,2).Value = "MO V1.4"
$ws.Cells.Item(7,3).Value = 42100
$ws.Cells.Item(7,5).Value = "dEFECT"
$ws.Cells.Item(7,6).Value = "MonV1.4: Qrcode: Medias are retrieved using media ID's as QR code"
$ws.Cells.Item(7,7).Value = "1. Choose Place media as the category`n2. Capture an image and save`n3. Select Qrcode option to identify media to link`n4. Enter 'media id' number in Qrcode generator`n5. Scan the Qrcode`nActual: The media retrieved has the 'media id' entered in the Qrcode generator`nExpected: The app must validate the QRcode and display a message informing that the qr code is not linked to any medias`nNote: When a media Id is given as a input in the Qrcode generator and scanned, the app must validate with the QRcodes present in the database and not with the media Id's"
$ws.Cells.Item(7,8).Value = "Android ICS"
$ws.Cells.Item(7,9).Value = "Functional"

# Match the (wrapped-text) row height Excel used for this new entry.
$ws.Rows.Item(7).RowHeight = 346.5

# --- Update the view/selection state to reflect the newly added row ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G7").Select()
